$d = $word.ActiveDocument

# 1. Merge "Specific heat " + "s" + "imulation" into a single run "Specific heat simulation"
#    (heading "Specific heat simulation")
$d.Content.Find.Execute(
    "Specific heat simulation", $true, $false, $false, $false, $false,
    $true, 1, $false, "Specific heat simulation", 2) | Out-Null

# 2. Add a trailing space after "...simply using" (before the inline equation)
$d.Content.Find.Execute(
    "graphical analysis; simply using", $true, $false, $false, $false, $false,
    $true, 1, $false, "graphical analysis; simply using ", 2) | Out-Null

# 3. Add a leading space before "on one (or even a handful)..." (after the inline equation)
$d.Content.Find.Execute(
    "on one (or even a handful) of measurements", $true, $false, $false, $false, $false,
    $true, 1, $false, " on one (or even a handful) of measurements", 2) | Out-Null

# 4. Merge "G" + "ive a brief description..." into a single run
$d.Content.Find.Execute(
    "Give a brief description of your process in the Materials & Methods section of your report; it doesn" + [char]8217 + "t need to be more than a couple of sentences.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Give a brief description of your process in the Materials & Methods section of your report; it doesn" + [char]8217 + "t need to be more than a couple of sentences.",
    2) | Out-Null
